$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.856

$ws.Range("B4").Value = 5.079
$ws.Range("C4").Value = -12.927
$ws.Range("E4").Value = 12.134

$ws.Range("C5").Value = -12.78

$ws.Range("B6").Value = 6.728
$ws.Range("C6").Value = -12.61

$ws.Range("B7").Value = 6.709000000000001

$ws.Range("B8").Value = 5.654999999999999
$ws.Range("C8").Value = -12.597

$ws.Range("E9").Value = 13.018

$ws.Range("E11").Value = 12.852

$ws.Range("E14").Value = 13.056

$ws.Range("B16").Value = 5.321999999999999
$ws.Range("C16").Value = -12.597

$ws.Range("E18").Value = 12.657

$ws.Range("B20").Value = 5.766999999999999

$ws.Range("B21").Value = 6.247

$ws.Range("C22").Value = -12.467

$ws.Range("E25").Value = 12.827
